$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Allocated Population") - shifts old Shelter
# Assigned/Shelter Level columns from B/C to C/D, carrying header style.
$ws.Columns.Item(2).Insert()
$ws.Cells.Item(1,2).Value2 = "Allocated Population"

# Community / Allocated Population / Shelter Assigned / Shelter Level
$data = @(
  @("Balite",       0,   "San Miguel Meysulao High School"),
  @("Balungao",      0,   "San Miguel Meysulao High School"),
  @("Bulusan",       0,   "San Miguel Meysulao High School"),
  @("Frances",       6,   "Frances E.C."),
  @("Gatbuca",       115, "Frances E.C."),
  @("Iba O'Este",    601, "San Miguel Meysulao High School")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value2 = $row[0]
  $ws.Cells.Item($r, 2).Value2 = $row[1]
  $ws.Cells.Item($r, 3).Value2 = $row[2]
  $ws.Cells.Item($r, 4).Value2 = 1
  $r = $r + 1
}
